$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2192987281846225
$ws.Range("D2").Value = 0.3777406617731509
$ws.Range("G2").Value = 0.1237476138499915
$ws.Range("H2").Value = 0.992
